$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.281766
$ws.Range("H2").Value = 45.845298
$ws.Range("I2").Value = 0.1817381432449346
$ws.Range("J2").Value = 0.1817381432449346
$ws.Range("M2").Value = 58.95713633333333
$ws.Range("N2").Value = 176.871409
$ws.Range("O2").Value = 0.4863146960083892
$ws.Range("P2").Value = 0.4863146960083893
$ws.Range("Q2").Value = 900.969161476098
$ws.Range("R2").Value = 8108.722453284882
$ws.Range("S2").Value = 0.08838192988528948
$ws.Range("T2").Value = 0.08838192988528948
$ws.Range("G3").Value = 15.281766
$ws.Range("H3").Value = 45.845298
$ws.Range("I3").Value = 0.1817381432449346
$ws.Range("J3").Value = 0.1817381432449346
$ws.Range("O3").Value = 0.07416766570679004
$ws.Range("P3").Value = 0.07416766570679005
$ws.Range("Q3").Value = 137.40645744198
$ws.Range("R3").Value = 1236.65811697782
$ws.Range("S3").Value = 0.01347909385436303
$ws.Range("T3").Value = 0.01347909385436303
$ws.Range("G4").Value = 15.281766
$ws.Range("H4").Value = 45.845298
$ws.Range("I4").Value = 0.1817381432449346
$ws.Range("J4").Value = 0.1817381432449346
$ws.Range("M4").Value = 42.51661933333333
$ws.Range("N4").Value = 127.549858
$ws.Range("O4").Value = 0.3507032073181665
$ws.Range("P4").Value = 0.3507032073181665
$ws.Range("Q4").Value = 649.7290277630759
$ws.Range("R4").Value = 5847.561249867684
$ws.Range("S4").Value = 0.06373614972804693
$ws.Range("T4").Value = 0.06373614972804695
$ws.Range("G5").Value = 15.281766
$ws.Range("H5").Value = 45.845298
$ws.Range("I5").Value = 0.1817381432449346
$ws.Range("J5").Value = 0.1817381432449346
$ws.Range("M5").Value = 10.76719366666667
$ws.Range("N5").Value = 32.301581
$ws.Range("O5").Value = 0.0888144309666542
$ws.Range("P5").Value = 0.08881443096665421
$ws.Range("Q5").Value = 164.541734090682
$ws.Range("R5").Value = 1480.875606816138
$ws.Range("S5").Value = 0.01614096977723516
$ws.Range("T5").Value = 0.01614096977723516
$ws.Range("I6").Value = 0.2947137116012682
$ws.Range("J6").Value = 0.2947137116012682
$ws.Range("M6").Value = 58.95713633333333
$ws.Range("N6").Value = 176.871409
$ws.Range("O6").Value = 0.4863146960083892
$ws.Range("P6").Value = 0.4863146960083893
$ws.Range("Q6").Value = 1461.046981530136
$ws.Range("R6").Value = 13149.42283377122
$ws.Range("S6").Value = 0.1433236090668749
$ws.Range("T6").Value = 0.1433236090668749
$ws.Range("I7").Value = 0.2947137116012682
$ws.Range("J7").Value = 0.2947137116012682
$ws.Range("O7").Value = 0.07416766570679004
$ws.Range("P7").Value = 0.07416766570679005
$ws.Range("S7").Value = 0.02185822804125019
$ws.Range("T7").Value = 0.02185822804125019
$ws.Range("I8").Value = 0.2947137116012682
$ws.Range("J8").Value = 0.2947137116012682
$ws.Range("M8").Value = 42.51661933333333
$ws.Range("N8").Value = 127.549858
$ws.Range("O8").Value = 0.3507032073181665
$ws.Range("P8").Value = 0.3507032073181665
$ws.Range("Q8").Value = 1053.626112208432
$ws.Range("R8").Value = 9482.635009875888
$ws.Range("S8").Value = 0.1033570438992059
$ws.Range("T8").Value = 0.1033570438992059
$ws.Range("I9").Value = 0.2947137116012682
$ws.Range("J9").Value = 0.2947137116012682
$ws.Range("M9").Value = 10.76719366666667
$ws.Range("N9").Value = 32.301581
$ws.Range("O9").Value = 0.0888144309666542
$ws.Range("P9").Value = 0.08881443096665421
$ws.Range("Q9").Value = 266.827339056824
$ws.Range("R9").Value = 2401.446051511416
$ws.Range("S9").Value = 0.02617483059393727
$ws.Range("T9").Value = 0.02617483059393728
$ws.Range("G10").Value = 18.371237
$ws.Range("H10").Value = 55.113711
$ws.Range("I10").Value = 0.2184796247693259
$ws.Range("J10").Value = 0.2184796247693259
$ws.Range("M10").Value = 58.95713633333333
$ws.Range("N10").Value = 176.871409
$ws.Range("O10").Value = 0.4863146960083892
$ws.Range("P10").Value = 0.4863146960083893
$ws.Range("Q10").Value = 1083.115524420978
$ws.Range("R10").Value = 9748.0397197888
$ws.Range("S10").Value = 0.1062498523037217
$ws.Range("T10").Value = 0.1062498523037217
$ws.Range("G11").Value = 18.371237
$ws.Range("H11").Value = 55.113711
$ws.Range("I11").Value = 0.2184796247693259
$ws.Range("J11").Value = 0.2184796247693259
$ws.Range("O11").Value = 0.07416766570679004
$ws.Range("P11").Value = 0.07416766570679005
$ws.Range("Q11").Value = 165.18552862261
$ws.Range("R11").Value = 1486.66975760349
$ws.Range("S11").Value = 0.01620412377363629
$ws.Range("T11").Value = 0.01620412377363629
$ws.Range("G12").Value = 18.371237
$ws.Range("H12").Value = 55.113711
$ws.Range("I12").Value = 0.2184796247693259
$ws.Range("J12").Value = 0.2184796247693259
$ws.Range("M12").Value = 42.51661933333333
$ws.Range("N12").Value = 127.549858
$ws.Range("O12").Value = 0.3507032073181665
$ws.Range("P12").Value = 0.3507032073181665
$ws.Range("Q12").Value = 781.0828902114487
$ws.Range("R12").Value = 7029.746011903038
$ws.Range("S12").Value = 0.07662150514027213
$ws.Range("T12").Value = 0.07662150514027215
$ws.Range("G13").Value = 18.371237
$ws.Range("H13").Value = 55.113711
$ws.Range("I13").Value = 0.2184796247693259
$ws.Range("J13").Value = 0.2184796247693259
$ws.Range("M13").Value = 10.76719366666667
$ws.Range("N13").Value = 32.301581
$ws.Range("O13").Value = 0.0888144309666542
$ws.Range("P13").Value = 0.08881443096665421
$ws.Range("Q13").Value = 197.8066666752323
$ws.Range("R13").Value = 1780.260000077091
$ws.Range("S13").Value = 0.01940414355169581
$ws.Range("T13").Value = 0.01940414355169582
$ws.Range("G14").Value = 25.652214
$ws.Range("H14").Value = 76.956642
$ws.Range("I14").Value = 0.3050685203844711
$ws.Range("J14").Value = 0.3050685203844711
$ws.Range("M14").Value = 58.95713633333333
$ws.Range("N14").Value = 176.871409
$ws.Range("O14").Value = 0.4863146960083892
$ws.Range("P14").Value = 0.4863146960083893
$ws.Range("Q14").Value = 1512.381078049842
$ws.Range("R14").Value = 13611.42970244858
$ws.Range("S14").Value = 0.1483593047525032
$ws.Range("T14").Value = 0.1483593047525032
$ws.Range("G15").Value = 25.652214
$ws.Range("H15").Value = 76.956642
$ws.Range("I15").Value = 0.3050685203844711
$ws.Range("J15").Value = 0.3050685203844711
$ws.Range("O15").Value = 0.07416766570679004
$ws.Range("P15").Value = 0.07416766570679005
$ws.Range("Q15").Value = 230.65265174742
$ws.Range("R15").Value = 2075.87386572678
$ws.Range("S15").Value = 0.02262622003754052
$ws.Range("T15").Value = 0.02262622003754052
$ws.Range("G16").Value = 25.652214
$ws.Range("H16").Value = 76.956642
$ws.Range("I16").Value = 0.3050685203844711
$ws.Range("J16").Value = 0.3050685203844711
$ws.Range("M16").Value = 42.51661933333333
$ws.Range("N16").Value = 127.549858
$ws.Range("O16").Value = 0.3507032073181665
$ws.Range("P16").Value = 0.3507032073181665
$ws.Range("Q16").Value = 1090.645417695204
$ws.Range("R16").Value = 9815.808759256837
$ws.Range("S16").Value = 0.1069885085506415
$ws.Range("T16").Value = 0.1069885085506415
$ws.Range("G17").Value = 25.652214
$ws.Range("H17").Value = 76.956642
$ws.Range("I17").Value = 0.3050685203844711
$ws.Range("J17").Value = 0.3050685203844711
$ws.Range("M17").Value = 10.76719366666667
$ws.Range("N17").Value = 32.301581
$ws.Range("O17").Value = 0.0888144309666542
$ws.Range("P17").Value = 0.08881443096665421
$ws.Range("Q17").Value = 276.202356116778
$ws.Range("R17").Value = 2485.821205051002
$ws.Range("S17").Value = 0.02709448704378595
$ws.Range("T17").Value = 0.02709448704378595
